# Adding open extent report feature
# - Adds a new row (DS002 / QA) to the Login sheet, mirroring the existing
#   DS001 row (same style/hyperlink pattern), and makes the Login sheet the
#   active sheet/selection instead of the Address sheet.

$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item("Login")

# Register the hyperlink relationship first (target cell gets overwritten
# with the correct style/value by the row copy below).
$loginSheet.Hyperlinks.Add($loginSheet.Range("C3"), "https://www.bigbasket.com/")

# Duplicate row 2 (style + values) down into row 3, exactly like the
# existing DS001 row, then overwrite the two cells that differ.
$loginSheet.Range("A2:F2").Copy($loginSheet.Range("A3:F3"))
$loginSheet.Range("A3").Value = "DS002"
$loginSheet.Range("B3").Value = "QA"

# Make the Login sheet the active tab/sheet with A3 selected, which also
# clears the "Address" sheet's tabSelected flag.
$loginSheet.Activate()
[void]$loginSheet.Range("A3").Select()
